$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell I2: "Approved" -> "Rejected"
$ws.Range("I2").Value = "Rejected"

# Add new cell J2 with value "Testing"
$ws.Range("J2").Value = "Testing"

# Update the active selection to H14
$ws.Range("H14").Select()
